# Word COM-interop script: append three new "to do" list items after the
# last paragraph of the document, moving the _GoBack bookmark into the
# newly added third item (between "dabei " and "waren").

$d = $word.ActiveDocument

# --- 1) Drop the _GoBack bookmark from the current last paragraph ---
# (Deleting the bookmark surgically keeps every other attribute of that
# paragraph - e.g. its w:rsidR/w:rsidP - untouched.)
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# --- 2) Append the three new list paragraphs after it ---
$lastPara = $d.Paragraphs.Last
$insertionPoint = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End - 1)

$newParasXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">Datensatz </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>gv_children</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">: Financial </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>gifts</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>to</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>children</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> näher betrachten </w:t></w:r><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> nützlich?</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">Location </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>child</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> aus Welle 5 übernehmen</w:t></w:r><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> (ggf. weiter zurück)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">Welle, in der sie das erste Mal dabei </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>waren</w:t></w:r></w:p>
'@

[void]$insertionPoint.InsertXML($newParasXml)
